# Automatic Excel update [2025-07-31 07:27:43]
#
# Bumps the "Data ostatniej aktualizacji" (E column) from 2025-07-30 to
# 2025-07-31 for every listing that was re-checked that day, on both
# sheets. Three listings flipped from Active to Inactive (H column,
# "Aktywne") instead of being re-verified, so their date stays put.

$wb = $excel.ActiveWorkbook

# Rows (1-based worksheet rows) whose E cell should move from
# 2025-07-30 -> 2025-07-31, per sheet.
$sheet1Rows = @(2,3,4,5,6,7,8,9,10,11,13,15,16,17,18,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,46,47,48,49,50,51,52,53,54,55,56,57,58,59,61,62,63,64)
$sheet2Rows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,16,18,19,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,41,42,43,44,45,46,47,48,49,50,51,52,53,55,56,58,59)

function Update-Dates {
    param($ws, $rows)

    foreach ($r in $rows) {
        $cell = $ws.Cells.Item($r, 5)   # column E
        if ($cell.Value2 -eq "2025-07-30") {
            # Leading apostrophe forces literal text so the engine doesn't
            # reinterpret the date-shaped string as a serial date number.
            $cell.Value = "'2025-07-31"
        }
    }
}

$ws1 = $wb.Worksheets.Item("powiat krakowski")
Update-Dates $ws1 $sheet1Rows
$ws1.Cells.Item(65, 8).Value = $false   # H65 "Aktywne" -> FALSE

$ws2 = $wb.Worksheets.Item("powiat wielicki")
Update-Dates $ws2 $sheet2Rows
$ws2.Cells.Item(20, 8).Value = $false   # H20 "Aktywne" -> FALSE
$ws2.Cells.Item(60, 8).Value = $false   # H60 "Aktywne" -> FALSE
